# Refreshed cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Coin, Link, Price, Volume(1h) -- written to columns B:E starting at row 2.
$cryptoData = @(
    @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "30.177.05", "  +0.37%  "),
    @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.911.41", "  -0.02%  "),
    @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.002", "  +0.09%  "),
    @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.8215", "  +4.44%  "),
    @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "244.14", "  +0.45%  "),
    @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.001", "  +0.08%  "),
    @("LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "1.909.64", "  +0.04%  "),
    @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3253", "  +2.73%  "),
    @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "26.86", "  +2.22%  "),
    @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.07051", "  +1.72%  "),
    @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.08104", "  +1.51%  "),
    @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.7727", "  +3.48%  "),
    @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.914.41", "  +0.09%  "),
    @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.295", "  +1.22%  "),
    @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "93.10", "  -0.36%  "),
    @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "30.174.19", "  +0.32%  "),
    @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "14.24", "  +1.52%  "),
    @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.913", "  -0.24%  "),
    @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "245.97", "  -0.55%  "),
    @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000007789", "  +0.04%  "),
    @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.165.44", "  -0.07%  "),
    @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.002", "  +0.18%  "),
    @("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.002", "  +0.13%  "),
    @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.063", "  +2.05%  "),
    @("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1657", "  +19.99%  "),
    @("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "9.333", "  +0.26%  "),
    @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "167.27", "  -1.61%  "),
    @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.99", "  +0.30%  "),
    @("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.110", "  +3.15%  "),
    @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.372", "  -0.19%  "),
    @("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.529", "  +0.38%  "),
    @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05855", "  +1.34%  "),
    @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.307", "  -0.83%  "),
    @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.106", "  -0.31%  "),
    @("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.275", "  +0.88%  "),
    @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7358", "  -0.18%  "),
    @("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.703", "  -0.71%  "),
    @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01926", "  +0.14%  "),
    @("MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.797", "  +0.09%  "),
    @("TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.4467", "  +0.52%  "),
    @("Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "73.48", "  +1.17%  "),
    @("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "5.963", "  -3.46%  "),
    @("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.8530", "  +2.21%  "),
    @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.912", "  +0.64%  "),
    @("PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.002", "  +0.10%  "),
    @("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "102.77", "  +2.16%  "),
    @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "9.918", "  +1.20%  "),
    @("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.594", "  +0.24%  "),
    @("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.005.97", "  +1.86%  "),
    @("RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "2.067.33", "  +0.38%  ")
)

$row = 2
foreach ($entry in $cryptoData) {
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $row = $row + 1
}

